# The deck currently has its "live" design (the theme wired to the slide
# master / used by every slide) holding the "Integral" color palette.
# The target commit swaps that for the default "Office Theme" palette
# (the complementary theme part that ships alongside it, previously
# unused/inert), i.e. the Design gallery selection was changed back to
# the stock Office theme.
#
# Reach the single reachable Theme/ColorScheme object (every slide
# shares the one Design/SlideMaster in this deck) and push each of the
# 12 theme colors to the standard Office palette. PowerPoint's
# ColorScheme.Colors(i).RGB uses COLORREF order (0xBBGGRR), so the RGB
# hex triples below are byte-swapped accordingly.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# index : element   : target RRGGBB -> COLORREF (0xBBGGRR)
$colors.Colors(1).RGB  = 0x000000   # dk1      000000
$colors.Colors(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$colors.Colors(3).RGB  = 0x6A5444   # dk2      44546A
$colors.Colors(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$colors.Colors(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$colors.Colors(6).RGB  = 0x317DED   # accent2  ED7D31
$colors.Colors(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$colors.Colors(8).RGB  = 0x00C0FF   # accent4  FFC000
$colors.Colors(9).RGB  = 0xC47244   # accent5  4472C4
$colors.Colors(10).RGB = 0x47AD70   # accent6  70AD47
$colors.Colors(11).RGB = 0xC16305   # hlink    0563C1
$colors.Colors(12).RGB = 0x724F95   # folHlink 954F72

# Mirror the naming too (no-op on hosts that don't serialize it, but
# matches the Design gallery rename PowerPoint performs when swapping
# back to a built-in theme).
$design.Name = "Office Theme"
$theme.Name = "Office Theme"
